# Updates the cryptocurrency prices/volumes (and, for a few rows, the
# coin name/link because the source ranking reordered) to the latest
# scraped values. Price strings that look like plain numbers (e.g.
# "231.51") are written with a leading apostrophe so Excel keeps them
# as text (matching the original "Price" column formatting) instead of
# auto-converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.825.46'
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").Value = '1.810.40'
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("E4").Value = '  +0.54%  '
$ws.Range("D5").Value = '''231.51'
$ws.Range("E5").Value = '  +3.08%  '
$ws.Range("D6").Value = '''0.608'
$ws.Range("E6").Value = '  +1.22%  '
$ws.Range("E7").Value = '  +0.47%  '
$ws.Range("D8").Value = '''39.70'
$ws.Range("E8").Value = '  -4.45%  '
$ws.Range("D9").Value = '''0.307'
$ws.Range("E9").Value = '  +5.21%  '
$ws.Range("E10").Value = '  +2.23%  '
$ws.Range("E11").Value = '  +0.24%  '
$ws.Range("D12").Value = '2.073.21'
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("D13").Value = '1.816.72'
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("D14").Value = '''11.04'
$ws.Range("E14").Value = '  +1.07%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '''4.64'
$ws.Range("E15").Value = '  +5.54%  '
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").Value = '''0.655'
$ws.Range("E16").Value = '  +4.26%  '
$ws.Range("D17").Value = '34.826.19'
$ws.Range("E17").Value = '  +0.96%  '
$ws.Range("D18").Value = '''68.70'
$ws.Range("E18").Value = '  +2.13%  '
$ws.Range("E19").Value = '  +1.83%  '
$ws.Range("D20").Value = '''236.74'
$ws.Range("E20").Value = '  -1.66%  '
$ws.Range("D21").Value = '''11.69'
$ws.Range("E21").Value = '  +4.76%  '
$ws.Range("D22").Value = '''4.70'
$ws.Range("E22").Value = '  +8.32%  '
$ws.Range("E23").Value = '  +0.46%  '
$ws.Range("E24").Value = '  +4.19%  '
$ws.Range("D25").Value = '''172.95'
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("D26").Value = '''7.72'
$ws.Range("E26").Value = '  +0.80%  '
$ws.Range("D27").Value = '''17.31'
$ws.Range("E27").Value = '  -0.45%  '
$ws.Range("D28").Value = '''0.120'
$ws.Range("E28").Value = '  -0.63%  '
$ws.Range("D29").Value = '''1.59'
$ws.Range("E29").Value = '  +30.23%  '
$ws.Range("E30").Value = '  +0.60%  '
$ws.Range("D31").Value = '3.339.77'
$ws.Range("E31").Value = '  +37.46%  '
$ws.Range("E32").Value = '  +6.18%  '
$ws.Range("E33").Value = '  +1.68%  '
$ws.Range("D34").Value = '''3.93'
$ws.Range("E34").Value = '  +1.43%  '
$ws.Range("D35").Value = '''1.78'
$ws.Range("E35").Value = '  -1.13%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").Value = '''1.18'
$ws.Range("E36").Value = '  +11.43%  '
$ws.Range("B37").Value = 'Aave'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D37").Value = '''93.43'
$ws.Range("E37").Value = '  +6.02%  '
$ws.Range("E38").Value = '  +4.18%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.0191'
$ws.Range("E39").Value = '  +1.95%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.304.29'
$ws.Range("E40").Value = '  -0.87%  '
$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").Value = '''1.27'
$ws.Range("E41").Value = '  +3.70%  '
$ws.Range("D42").Value = '''14.83'
$ws.Range("E42").Value = '  +0.38%  '
$ws.Range("E43").Value = '  +4.66%  '
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("E45").Value = '  +0.44%  '
$ws.Range("E46").Value = '  -1.27%  '
$ws.Range("D47").Value = '''6.19'
$ws.Range("E47").Value = '  +7.16%  '
$ws.Range("D48").Value = '''0.0512'
$ws.Range("E48").Value = '  -1.27%  '
$ws.Range("D49").Value = '1.988.20'
$ws.Range("E49").Value = '  +1.33%  '
$ws.Range("E50").Value = '  +0.42%  '
$ws.Range("D51").Value = '''0.0644'
$ws.Range("E51").Value = '  +5.89%  '
